# Generate Report for Handoff
# Adds a new row (row 3) to each sheet (Overview, zh-cn, de-de) for the
# newly-ready-for-handoff file c1710a5e-758d-4da5-9c23-9b911e2920cb.md

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet "Overview"
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Overview")

$ws.Range("A3").Value = "c1710a5e-758d-4da5-9c23-9b911e2920cb.md"
$ws.Range("B3").Value = "Ready for handoff"
$ws.Range("C3").Value = "Ready for handoff"
$ws.Range("D3").Value = "2016-03-24 11:28:45"
$ws.Range("D3").NumberFormat = "yyyy-mm-dd HH:mm:ss"

$ws.Hyperlinks.Add($ws.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/47ee042c2a70cc8972a9113b4b07e01fed7c4226/e2e/c1710a5e-758d-4da5-9c23-9b911e2920cb.md", "", "", "c1710a5e-758d-4da5-9c23-9b911e2920cb.md")

# ---------------------------------------------------------------------
# Sheet "zh-cn"
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("zh-cn")

$ws.Range("A3").Value = "c1710a5e-758d-4da5-9c23-9b911e2920cb.md"
$ws.Range("B3").Value = ".md"
$ws.Range("C3").Value = "Ready for handoff"
$ws.Range("D3").Value = "c1710a5e-758d-4da5-9c23-9b911e2920cb.42c1990588a1ca31ed4474337322df42ffb845e8.zh-cn.xlf"
$ws.Range("E3").Value = "2016-03-24 11:28:40"
$ws.Range("E3").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$ws.Range("H3").Value = "0001-01-01 00:00:00"
$ws.Range("H3").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$ws.Range("J3").Value = "Include"

$ws.Hyperlinks.Add($ws.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/47ee042c2a70cc8972a9113b4b07e01fed7c4226/e2e/c1710a5e-758d-4da5-9c23-9b911e2920cb.md", "", "", "c1710a5e-758d-4da5-9c23-9b911e2920cb.md")
$ws.Hyperlinks.Add($ws.Range("D3"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/b10663b1e1bd3d8193a32f46b67e51046385190e/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/c1710a5e-758d-4da5-9c23-9b911e2920cb.42c1990588a1ca31ed4474337322df42ffb845e8.zh-cn.xlf", "", "", "c1710a5e-758d-4da5-9c23-9b911e2920cb.42c1990588a1ca31ed4474337322df42ffb845e8.zh-cn.xlf")

# ---------------------------------------------------------------------
# Sheet "de-de"
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("de-de")

$ws.Range("A3").Value = "c1710a5e-758d-4da5-9c23-9b911e2920cb.md"
$ws.Range("B3").Value = ".md"
$ws.Range("C3").Value = "Ready for handoff"
$ws.Range("D3").Value = "c1710a5e-758d-4da5-9c23-9b911e2920cb.42c1990588a1ca31ed4474337322df42ffb845e8.de-de.xlf"
$ws.Range("E3").Value = "2016-03-24 11:28:45"
$ws.Range("E3").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$ws.Range("H3").Value = "0001-01-01 00:00:00"
$ws.Range("H3").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$ws.Range("J3").Value = "Include"

$ws.Hyperlinks.Add($ws.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/47ee042c2a70cc8972a9113b4b07e01fed7c4226/e2e/c1710a5e-758d-4da5-9c23-9b911e2920cb.md", "", "", "c1710a5e-758d-4da5-9c23-9b911e2920cb.md")
$ws.Hyperlinks.Add($ws.Range("D3"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/95b6e116595c8776948c99994be0755bd3d9ea52/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/c1710a5e-758d-4da5-9c23-9b911e2920cb.42c1990588a1ca31ed4474337322df42ffb845e8.de-de.xlf", "", "", "c1710a5e-758d-4da5-9c23-9b911e2920cb.42c1990588a1ca31ed4474337322df42ffb845e8.de-de.xlf")

Write-Output "Handoff report row added to Overview, zh-cn, and de-de sheets."
